$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextCell "D2" "65.857.41"
Set-TextCell "E2" "  +1.74%  "
Set-TextCell "D3" "2.695.63"
Set-TextCell "E3" "  +2.40%  "
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "609.19"
Set-TextCell "E5" "  +2.81%  "
Set-TextCell "D6" "158.24"
Set-TextCell "E6" "  +2.29%  "
Set-TextCell "E7" "  -0.03%  "
Set-TextCell "E8" "  -0.50%  "
Set-TextCell "E9" "  +6.10%  "
Set-TextCell "D10" "6.06"
Set-TextCell "E10" "  +5.31%  "
Set-TextCell "E11" "  +1.70%  "
Set-TextCell "E12" "  +1.23%  "
Set-TextCell "D13" "30.24"
Set-TextCell "E13" "  +4.91%  "
Set-TextCell "E14" "  +9.79%  "
Set-TextCell "D15" "3.182.35"
Set-TextCell "E15" "  +2.55%  "
Set-TextCell "D16" "65.734.06"
Set-TextCell "E16" "  +1.60%  "
Set-TextCell "D17" "2.701.05"
Set-TextCell "E17" "  +4.31%  "
Set-TextCell "D18" "12.66"
Set-TextCell "E18" "  +1.18%  "
Set-TextCell "E19" "  +2.62%  "
Set-TextCell "D20" "359.88"
Set-TextCell "E20" "  +2.84%  "
Set-TextCell "D21" "7.54"
Set-TextCell "E21" "  +4.15%  "
Set-TextCell "E22" "  -0.14%  "
Set-TextCell "D23" "70.65"
Set-TextCell "E23" "  +4.08%  "
Set-TextCell "D24" "9.86"
Set-TextCell "E24" "  +4.44%  "
Set-TextCell "D25" "0.0000106"
Set-TextCell "E25" "  +14.67%  "
Set-TextCell "E26" "  -2.43%  "
Set-TextCell "E27" "  +3.32%  "
Set-TextCell "E28" "  +5.79%  "
Set-TextCell "D29" "8.37"
Set-TextCell "E29" "  +3.71%  "
Set-TextCell "B30" "PancakeSwap"
Set-TextCell "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D30" "2.20"
Set-TextCell "E30" "  +5.83%  "
Set-TextCell "B31" "Bittensor"
Set-TextCell "C31" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D31" "544.26"
Set-TextCell "E31" "  +6.78%  "
Set-TextCell "E32" "  -1.55%  "
Set-TextCell "E33" "  +1.51%  "
Set-TextCell "D34" "6.73"
Set-TextCell "E34" "  +8.71%  "
Set-TextCell "D35" "5.38"
Set-TextCell "E35" "  -3.71%  "
Set-TextCell "D36" "0.433"
Set-TextCell "E36" "  +2.51%  "
Set-TextCell "D37" "20.81"
Set-TextCell "E37" "  +3.77%  "
Set-TextCell "D38" "163.14"
Set-TextCell "E38" "  -0.93%  "
Set-TextCell "E39" "  +0.00%  "
Set-TextCell "D40" "1.00"
Set-TextCell "E40" "  -0.03%  "
Set-TextCell "E41" "  -0.03%  "
Set-TextCell "B42" "Aave"
Set-TextCell "C42" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D42" "169.81"
Set-TextCell "E42" "  +3.71%  "
Set-TextCell "B43" "OKB"
Set-TextCell "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D43" "42.90"
Set-TextCell "E43" "  +1.65%  "
Set-TextCell "D44" "4.19"
Set-TextCell "E44" "  +3.41%  "
Set-TextCell "D45" "0.0616"
Set-TextCell "E45" "  +1.17%  "
Set-TextCell "D46" "23.66"
Set-TextCell "E46" "  +4.23%  "
Set-TextCell "E47" "  +4.78%  "
Set-TextCell "D48" "0.0267"
Set-TextCell "E48" "  +5.65%  "
Set-TextCell "E49" "  +2.35%  "
Set-TextCell "D50" "21.07"
Set-TextCell "E50" "  +9.61%  "
Set-TextCell "E51" "  +1.78%  "
